$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 22 de Marzo de 2020 a las 04:16"

# Swap the Huesca/Huelva rows (A53/A54) and their "Casos activos" (C53/C54) values
$ws.Range("A53").Value = "Huelva"
$ws.Range("A54").Value = "Huesca"

$ws.Range("C53").Value = 72
$ws.Range("C54").Value = 0
